$d = $word.ActiveDocument

# Helper: set a table cell's text robustly even when the cell's paragraph
# originally contained more than one run (Word's Range.Text setter here
# only overwrites the first run it touches, so trim any stale leftover
# text from runs further into the cell).
function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText

    $cell2 = $table.Cell($row, $col)
    $actual = $cell2.Range.Text
    if ($actual.Length -gt $newText.Length) {
        $rng = $cell2.Range
        $trimStart = $rng.Start + $newText.Length
        $trimEnd = $rng.End - 1
        if ($trimEnd -gt $trimStart) {
            $extra = $d.Range($trimStart, $trimEnd)
            $extra.Text = ""
        }
    }
}

# ---------------------------------------------------------------------
# Table 2: "Work item assignments"
# ---------------------------------------------------------------------
$workItems = $d.Tables.Item(2)

# Row 2 (Work item 1 - Submitting points of interest): Hours Worked 0 -> 5
Set-CellText $workItems 2 7 "5"

# Row 3 (Work item 2 - Displaying points of interest): State Incomplete -> Complete,
# Hours Worked 0 -> 4
Set-CellText $workItems 3 4 "Complete"
Set-CellText $workItems 3 7 "4"

# Row 4 (Work item 3 - Writing reviews): State "Inc"+"omplete" -> Complete,
# Hours Worked 0 -> 10
Set-CellText $workItems 4 4 "Complete"
Set-CellText $workItems 4 7 "10"

# Row 5 (Work item 4 - Pulling reviews and images): State Incomplete -> Complete,
# Hours Worked 0 -> 2
Set-CellText $workItems 5 4 "Complete"
Set-CellText $workItems 5 7 "2"

# ---------------------------------------------------------------------
# Table 3: "Issues"
# ---------------------------------------------------------------------
$issues = $d.Tables.Item(3)
Set-CellText $issues 2 1 "Networking library could not upload images"
Set-CellText $issues 2 2 "Resolved"
Set-CellText $issues 2 3 "The networking library, Volley, does not support multipart requests and therefore could not upload images as needed. This was resolved by switching to Retrofit."

# ---------------------------------------------------------------------
# Table 4: "Assessment"
# ---------------------------------------------------------------------
$assessment = $d.Tables.Item(4)
Set-CellText $assessment 1 2 "Entire iteration"
Set-CellText $assessment 2 2 "2020-07-27"
Set-CellText $assessment 3 2 "Tyler, Joey, Sanjay"
Set-CellText $assessment 4 2 "Green"

# ---------------------------------------------------------------------
# Final summary paragraph: append narrative text after the line break
# ---------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count - 1
$summaryPara = $d.Paragraphs.Item($lastParaIndex)
$summaryRange = $summaryPara.Range
$summaryRange.InsertAfter("Most items outlined in this iteration were completed, and the last one uncompleted will be completed within a day of the assessment. Sanjay was unable to complete his work due to illness, so Tyler did it instead. The app can now pull and display points of interest, and show and upload reviews and pictures for paths and points of interest. Submitting points of interest is almost done and will be completed within the next day.")
